# Update data provider for createAccountTest: extend AccountCreationData
# sheet with two more rows (3 and 4), re-using row 2 as a template, and
# give each row its own unique e-mail address plus mailto hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

# --- Duplicate row 2's formatting into rows 3 and 4 ----------------------
$ws.Range("A2:O2").Copy()
$ws.Range("A3:O3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:O2").Copy()
$ws.Range("A4:O4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Copy the (non-email) data values from row 2 into rows 3 and 4 -------
for ($c = 2; $c -le 15; $c++) {
    $v = $ws.Cells.Item(2, $c).Value2
    $ws.Cells.Item(3, $c).Value2 = $v
    $ws.Cells.Item(4, $c).Value2 = $v
}

# --- Remove the existing hyperlink (keeps cell text for now) -------------
$ws.Hyperlinks.Delete()

# --- Set the three unique e-mail addresses used by the data provider -----
$ws.Cells.Item(2, 1).Value2 = "newtest1@gmail.com"
$ws.Cells.Item(3, 1).Value2 = "newtest2@gmail.com"
$ws.Cells.Item(4, 1).Value2 = "newtest3@gmail.com"

# --- Re-create the mailto hyperlinks on A2, A3 and A4 ---------------------
$ws.Hyperlinks.Add($ws.Cells.Item(2, 1), "mailto:newtest1@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "mailto:newtest2@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 1), "mailto:newtest3@gmail.com") | Out-Null

# --- Adding hyperlinks re-styles column A, so restore the original look --
$ws.Range("A2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Select B4 like the authored workbook does ----------------------------
$ws.Range("B4").Select() | Out-Null
